$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K7").Value = -0.3352267436446591
$ws.Range("J8").Value = 0.01855976243503714
$ws.Range("I9").Value = -0.1296176279974082
$ws.Range("H10").Value = -0.2870636170015632
$ws.Range("G11").Value = 0.2135958395245076
$ws.Range("F12").Value = -0.06676204101096155
$ws.Range("E13").Value = 0.1052128168340501
$ws.Range("D14").Value = -0.2006497229122814
$ws.Range("C15").Value = 0.4116802297750048
$ws.Range("B16").Value = -0.2766911554241067
